$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new customer row 11: phone 71076781 (stored as text, like the
# source data), blank birthday, 0 points.
#
# A leading apostrophe forces Excel to store the entered value as text
# instead of a number (matching A11's t="inlineStr"/text type in the
# target). Excel normally flags such "number stored as text" cells with a
# quotePrefix style; resetting the style back to Normal afterwards drops
# that flag so the cell ends up with the plain/default style, same as the
# surrounding data cells.
$ws.Cells.Item(11, 1).Value = "'71076781"
$ws.Cells.Item(11, 1).Style = "Normal"

# B11 mirrors B10: present in the sheet but holding an empty string (not
# just "no cell"). A bare apostrophe enters an empty, text-typed value;
# clear the resulting quotePrefix style the same way as above.
$ws.Cells.Item(11, 2).Value = "'"
$ws.Cells.Item(11, 2).Style = "Normal"

$ws.Cells.Item(11, 3).Value = 0
